# Edit script: apply the StorageComponentClassDiagram.pptx changes
#   1. Update the cached "datetimeFigureOut" footer field text from
#      "3/12/2019" to "22-Mar-19" on the Slide Master and every Slide
#      Layout (the Header & Footer "Update automatically" date format
#      was changed, which re-renders the cached field text).
#   2. Rename "AddressBook" -> "TravelBuddy" everywhere it appears
#      (and the derived "AddressBookStorage" / "JsonAddressBook")
#      inside the class-diagram shapes on slide 1.

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes, $newDate) {
    $cnt = $shapes.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            $trLen = $tr.Length
            if ($trLen -gt 0) {
                $full = $tr.Characters(1, $trLen)
                $full.Text = $newDate
            }
        }
    }
}

$newDate = "22-Mar-19"

# --- Slide Master ---
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes $newDate

# --- Every Slide Layout belonging to the master ---
$layoutCount = $master.CustomLayouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes $newDate
}

# --- Rename AddressBook -> TravelBuddy on slide 1 ---
$s = $p.Slides.Item(1)

function Replace-SubText($shape, $start, $length, $newText) {
    $tr = $shape.TextFrame.TextRange
    $sub = $tr.Characters($start, $length)
    $sub.Text = $newText
}

$shapeCount = $s.Shapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $full = $shp.TextFrame.TextRange.Text
    $shapeId = $shp.Id

    if ($shapeId -eq 2) {
        # "<<interface>>" + line-break + "AddressBookStorage"
        $oldWord = "AddressBookStorage"
        $newWord = "TravelBuddyStorage"
        $wordLen = $oldWord.Length
        $idx = $full.IndexOf($oldWord)
        if ($idx -ge 0) {
            $startPos = $idx + 1
            Replace-SubText $shp $startPos $wordLen $newWord
        }
    }
    elseif ($shapeId -eq 50) {
        # "JsonAddressBook" + line-break + "Storage"
        $oldWord = "JsonAddressBook"
        $newWord = "JsonTravelBuddy"
        $wordLen = $oldWord.Length
        $idx = $full.IndexOf($oldWord)
        if ($idx -ge 0) {
            $startPos = $idx + 1
            Replace-SubText $shp $startPos $wordLen $newWord
        }
    }
    elseif ($shapeId -eq 66) {
        # "JsonSerializable" + line-break + "AddressBook"
        $oldWord = "AddressBook"
        $newWord = "TravelBuddy"
        $wordLen = $oldWord.Length
        $idx = $full.IndexOf($oldWord)
        if ($idx -ge 0) {
            $startPos = $idx + 1
            Replace-SubText $shp $startPos $wordLen $newWord
        }
    }
}
